# "ajustement locaux et staff"
# Add two new entries (Valcourt security, volunteers/CO showers) to the
# ADMIN (locaux) sheet, and move the active tab/selection from EQUIPES
# to ADMIN, positioned at the first empty row below the new data.

$wb = $excel.ActiveWorkbook
$wsAdmin = $wb.Worksheets.Item("ADMIN")

# New row 16: Valcourt security
$wsAdmin.Range("A16").Value = "Sécurité Valcourt"
$wsAdmin.Range("B16").Value = "Valcourt Securtiy services"
$wsAdmin.Range("C16").Value = 3110
$wsAdmin.Range("D16").Formula = "=C16"

# New row 17: Volunteers / CO showers
# (shared-string table is append-on-first-use, so write in the same
# order the original authoring tool did: French label, room code, then
# the English label)
$wsAdmin.Range("B17").Value = "Douches pour CO"
$wsAdmin.Range("C17").Value = "0240<br/>0250"
$wsAdmin.Range("D17").Formula = "=C17"
$wsAdmin.Range("A17").Value = "Volunteers showers"

# Make ADMIN the active sheet and move the selection to the next blank row
$wsAdmin.Activate()
$wsAdmin.Range("A18").Select()
